$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1 (00:16 -> 00:46)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 00:46"

# Swap the Huesca (row 53) and Huelva (row 54) rows so that Huelva comes
# first (row 53) with its "Casos activos" value of 72, and Huesca moves
# to row 54 with "Casos activos" of 0 (all other columns are identical
# between the two rows, B=37, D=37, E=0).
$ws.Range("A53").Value = "Huelva"
$ws.Range("B53").Value = 37
$ws.Range("C53").Value = 72
$ws.Range("D53").Value = 37
$ws.Range("E53").Value = 0

$ws.Range("A54").Value = "Huesca"
$ws.Range("B54").Value = 37
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 37
$ws.Range("E54").Value = 0
